# Fruta / hortaliza, semanal
# A new weekly record needs to be inserted at the top of the "Feria Lagunitas
# de Puerto Montt - Mango" data block (row 33), pushing the existing rows
# 33-82 down to 34-83 (dimension grows from A1:T82 to A1:T83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 33; this shifts rows 33-82 down to
# 34-83 and grows the sheet dimension automatically.
$ws.Rows.Item(33).Insert()

# The fixed columns (Mercado ID/Mercado/Region/Codreg/Tipo/Producto.../
# Variedad/Calidad/Unidad/Origen/Kg per unidad) repeat for every record in
# this block, so copy them straight down from the row that is now 34 (the
# record that used to be row 33).
$ws.Cells.Item(33, 1).Value2  = $ws.Cells.Item(34, 1).Value2   # Mercado ID
$ws.Cells.Item(33, 2).Value2  = $ws.Cells.Item(34, 2).Value2   # Mercado
$ws.Cells.Item(33, 3).Value2  = $ws.Cells.Item(34, 3).Value2   # Region
$ws.Cells.Item(33, 5).Value2  = $ws.Cells.Item(34, 5).Value2   # Codreg
$ws.Cells.Item(33, 6).Value2  = $ws.Cells.Item(34, 6).Value2   # Tipo
$ws.Cells.Item(33, 7).Value2  = $ws.Cells.Item(34, 7).Value2   # Producto ID
$ws.Cells.Item(33, 8).Value2  = $ws.Cells.Item(34, 8).Value2   # Producto
$ws.Cells.Item(33, 9).Value2  = $ws.Cells.Item(34, 9).Value2   # Categoria ID
$ws.Cells.Item(33, 10).Value2 = $ws.Cells.Item(34, 10).Value2  # Categoria
$ws.Cells.Item(33, 11).Value2 = $ws.Cells.Item(34, 11).Value2  # Variedad
$ws.Cells.Item(33, 12).Value2 = $ws.Cells.Item(34, 12).Value2  # Calidad
$ws.Cells.Item(33, 17).Value2 = $ws.Cells.Item(34, 17).Value2  # Unidad de comercializacion
$ws.Cells.Item(33, 18).Value2 = $ws.Cells.Item(34, 18).Value2  # Origen
$ws.Cells.Item(33, 20).Value2 = $ws.Cells.Item(34, 20).Value2  # Kg / unidad

# New weekly observation's own data.
$ws.Cells.Item(33, 4).Value2  = 44477   # Fecha
$ws.Cells.Item(33, 13).Value2 = 160     # Volumen
$ws.Cells.Item(33, 14).Value2 = 8000    # Precio minimo
$ws.Cells.Item(33, 15).Value2 = 9000    # Precio maximo
$ws.Cells.Item(33, 16).Value2 = 8500    # Precio promedio ponderado
$ws.Cells.Item(33, 19).Value2 = 2125    # Precio $/Kg
